$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.546.59'
$ws.Range('E2').Value = '  +6.16%  '
$ws.Range('D3').Value = '2.477.96'
$ws.Range('E3').Value = '  +7.56%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '570.89'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +5.47%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.54'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +11.50%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +4.29%  '
$ws.Range('D9').Value = '2.474.89'
$ws.Range('E9').Value = '  +7.53%  '
$ws.Range('E10').Value = '  +6.39%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.76'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +4.29%  '
$ws.Range('E12').Value = '  +1.02%  '
$ws.Range('E13').Value = '  +7.49%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.53'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +14.75%  '
$ws.Range('D15').Value = '2.917.12'
$ws.Range('E15').Value = '  +7.58%  '
$ws.Range('D16').Value = '63.383.37'
$ws.Range('E16').Value = '  +6.05%  '
$ws.Range('E17').Value = '  +10.45%  '
$ws.Range('D18').Value = '2.472.60'
$ws.Range('E18').Value = '  +6.98%  '
$ws.Range('E19').Value = '  +9.24%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '345.87'
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.33'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +8.14%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.87'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +6.34%  '
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.87'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +4.04%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.176'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +4.24%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  +14.17%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.27'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +7.02%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.32'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +12.08%  '
$ws.Range('D30').Value = '0.0₃0821'
$ws.Range('E30').Value = '  +15.31%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.89'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +18.90%  '
$ws.Range('E32').Value = '  +8.58%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '175.42'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.97%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.52'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +12.14%  '
$ws.Range('E35').Value = '  +6.25%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '19.08'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +7.92%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '372.45'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +20.04%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.52'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +11.69%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.72'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +15.06%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '40.32'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +6.24%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '151.75'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +11.30%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.75'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +10.14%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '20.85'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +12.82%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.603'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +6.90%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0968'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.60%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0527'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +7.83%  '
$ws.Range('D49').Value = '0.0₆0237'
$ws.Range('E49').Value = '  +6.01%  '
$ws.Range('E50').Value = '  +7.33%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '18.30'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +10.09%  '

Write-Host "Applied all crypto list updates"
